# Updates the cryptos price list (row-by-row refresh of Price/Volume(1h),
# plus a rank swap between Hedera and EthereumClassic in rows 32-33).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.818.00'
$ws.Range('E2').Value = '  +0.61%  '
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.303.61'
$ws.Range('E3').Value = '  +0.63%  '
# Row 4
$ws.Range('E4').Value = '  +0.29%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '115.56'
$ws.Range('E5').Value = '  +21.36%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.10'
$ws.Range('E6').Value = '  +0.55%  '
# Row 7
$ws.Range('E7').Value = '  +0.57%  '
# Row 8
$ws.Range('E8').Value = '  +0.22%  '
# Row 9
$ws.Range('E9').Value = '  +2.76%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '48.87'
$ws.Range('E10').Value = '  +9.74%  '
# Row 11
$ws.Range('E11').Value = '  +0.74%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.71'
$ws.Range('E12').Value = '  +11.72%  '
# Row 13
$ws.Range('E13').Value = '  +1.99%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.69'
$ws.Range('E14').Value = '  +3.45%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.650.86'
$ws.Range('E15').Value = '  +0.76%  '
# Row 16
$ws.Range('E16').Value = '  +1.64%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.303.82'
$ws.Range('E17').Value = '  +0.58%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.697.49'
$ws.Range('E18').Value = '  +0.27%  '
# Row 19
$ws.Range('E19').Value = '  +3.06%  '
# Row 20
$ws.Range('E20').Value = '  +6.60%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.85'
$ws.Range('E21').Value = '  +0.47%  '
# Row 22
$ws.Range('E22').Value = '  +4.37%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.87'
$ws.Range('E23').Value = '  -0.05%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.72'
$ws.Range('E24').Value = '  +7.69%  '
# Row 25
$ws.Range('E25').Value = '  +15.17%  '
# Row 26
$ws.Range('E26').Value = '  -0.04%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.54'
$ws.Range('E27').Value = '  +3.02%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '44.37'
$ws.Range('E28').Value = '  +9.41%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.41'
$ws.Range('E29').Value = '  -1.46%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.27'
$ws.Range('E30').Value = '  -0.34%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '177.71'
$ws.Range('E31').Value = '  +1.39%  '
# Row 32
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0942'
$ws.Range('E32').Value = '  +6.94%  '
# Row 33
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.87'
$ws.Range('E33').Value = '  -0.25%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.56'
$ws.Range('E34').Value = '  +3.88%  '
# Row 35
$ws.Range('E35').Value = '  +1.08%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.75'
$ws.Range('E36').Value = '  +8.30%  '
# Row 37
$ws.Range('E37').Value = '  +2.08%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.95'
$ws.Range('E38').Value = '  +19.43%  '
# Row 39
$ws.Range('E39').Value = '  +0.28%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '74.80'
$ws.Range('E40').Value = '  +15.60%  '
# Row 41
$ws.Range('E41').Value = '  +3.96%  '
# Row 42
$ws.Range('E42').Value = '  +3.53%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.50'
$ws.Range('E43').Value = '  +12.41%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.44'
$ws.Range('E44').Value = '  +7.61%  '
# Row 45
$ws.Range('E45').Value = '  +0.03%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.95'
$ws.Range('E46').Value = '  +13.86%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.84'
$ws.Range('E47').Value = '  +0.23%  '
# Row 48
$ws.Range('E48').Value = '  -0.86%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '101.99'
$ws.Range('E49').Value = '  +3.96%  '
# Row 50
$ws.Range('E50').Value = '  +4.48%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.456'
$ws.Range('E51').Value = '  +6.85%  '
